# namapovani_poli.xlsx -- "czc xml scheme added"
#
# Renames the sole sheet to "Alza" and adds two new field-mapping sheets,
# "CZC" and "Onlineshop", each carrying two header rows: row 1 holds the
# internal/Helios field names (reusing the ones already used on the Alza
# sheet where they match), row 2 holds the external XML/eshop field names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing sheet, add the two new sheets right after it.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Alza"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "CZC"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Onlineshop"

# ---------------------------------------------------------------------
# 2. Carry over the two header-row styles already used on "Alza":
#    - A1/A2-style: bold text + thin border (file-name / first column)
#    - B1..-style : thin border only
#    Copy-format (not copy-value) so the existing style entries get
#    reused instead of Excel fabricating near-duplicate ones.
# ---------------------------------------------------------------------
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("A2").PasteSpecial(-4122)
$ws3.Range("A1").PasteSpecial(-4122)
$ws3.Range("A2").PasteSpecial(-4122)

$ws1.Range("B1").Copy()
$ws2.Range("B1:S1").PasteSpecial(-4122)
$ws2.Range("B2:S2").PasteSpecial(-4122)
$ws3.Range("B1:K1").PasteSpecial(-4122)
$ws3.Range("B2:K2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. CZC sheet content.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = 'sklad_helios.xlsx'
$ws2.Range("B1").Value = 'Registrační číslo'
$ws2.Range("C1").Value = 'Registrační číslo'
$ws2.Range("D1").Value = 'Název 1'
$ws2.Range("E1").Value = 'Čárový kód'
$ws2.Range("F1").Value = 'Prodejní cena'
$ws2.Range("G1").Value = 'Množství skladem'
$ws2.Range("H1").Value = 'MANUFACTURER'
$ws2.Range("I1").Value = 'IMAGE'
$ws2.Range("J1").Value = 'URL'
$ws2.Range("K1").Value = 'DESCRIPTION'
$ws2.Range("L1").Value = 'Currency'
$ws2.Range("M1").Value = 'CopyrightFee'
$ws2.Range("N1").Value = 'JC hist. recykl. přísp.'
$ws2.Range("O1").Value = 'Hmotnost'
$ws2.Range("P1").Value = 'Šířka'
$ws2.Range("Q1").Value = 'Výška'
$ws2.Range("R1").Value = 'Hloubka'
$ws2.Range("S1").Value = 'WARRANTY'

$ws2.Range("A2").Value = 'CZC_HELIOS.xml'
$ws2.Range("B2").Value = 'CODE'
$ws2.Range("C2").Value = 'ESHOP_CODE'
$ws2.Range("D2").Value = 'NAME'
$ws2.Range("E2").Value = 'EAN'
$ws2.Range("F2").Value = 'PRICE'
$ws2.Range("G2").Value = 'QUANTITY'
$ws2.Range("H2").Value = 'MANUFACTURER'
$ws2.Range("I2").Value = 'IMAGE'
$ws2.Range("J2").Value = 'URL'
$ws2.Range("K2").Value = 'DESCRIPTION'
$ws2.Range("L2").Value = 'CURRENCY'
$ws2.Range("M2").Value = 'AUTHOR_FEE'
$ws2.Range("N2").Value = 'RECYCLE_FEE'
$ws2.Range("O2").Value = 'WEIGHT_BRUTTO'
$ws2.Range("P2").Value = 'SIZE_X_NETTO'
$ws2.Range("Q2").Value = 'SIZE_Y_NETTO'
$ws2.Range("R2").Value = 'SIZE_Z_NETTO'
$ws2.Range("S2").Value = 'WARRANTY'

# ---------------------------------------------------------------------
# 4. Onlineshop sheet content.
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = 'sklad_helios.xlsx'
$ws3.Range("B1").Value = 'Registrační číslo'
$ws3.Range("C1").Value = 'Čárový kód'
$ws3.Range("D1").Value = 'Název 1'
$ws3.Range("E1").Value = 'MANUFACTURER'
$ws3.Range("F1").Value = 'DESCRIPTION'
$ws3.Range("G1").Value = 'Množství skladem'
$ws3.Range("H1").Value = 'Hmotnost'
$ws3.Range("I1").Value = 'IMAGE'
$ws3.Range("J1").Value = 'Prodejní cena'
$ws3.Range("K1").Value = 'Cena v HM'

$ws3.Range("A2").Value = 'onlineshop.xml'
$ws3.Range("B2").Value = 'ITEM_ID'
$ws3.Range("C2").Value = 'EAN'
$ws3.Range("D2").Value = 'PRODUCT'
$ws3.Range("E2").Value = 'MANUFACTURER'
$ws3.Range("F2").Value = 'DESCRIPTION'
$ws3.Range("G2").Value = 'STOCK'
$ws3.Range("H2").Value = 'WEIGHT'
$ws3.Range("I2").Value = 'IMGURL'
$ws3.Range("J2").Value = 'PRICE_NAKUP'
$ws3.Range("K2").Value = 'PRICE_DOPORUCENA'

# ---------------------------------------------------------------------
# 5. Active sheet / selection to match the final document (CZC active).
# ---------------------------------------------------------------------
$ws2.Select()
$ws2.Range("S1").Select()
